$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5859.8
$ws.Range("I18").Value = 10049.5
$ws.Range("J18").Value = 3066.6667
$ws.Range("K18").Value = 10049.5
$ws.Range("L18").Value = 3066.6667
$ws.Range("M18").Value = -9765.5
$ws.Range("N18").Value = -3634.6667
$ws.Range("H58").Value = 3453
$ws.Range("I58").Value = 671.6667
$ws.Range("J58").Value = 7625
$ws.Range("K58").Value = 2015.0001
$ws.Range("L58").Value = 22875
$ws.Range("M58").Value = -1865.0001
$ws.Range("N58").Value = -23175
$ws.Range("H86").Value = 5473.4
$ws.Range("I86").Value = 4823
$ws.Range("J86").Value = 6449
$ws.Range("K86").Value = 4823
$ws.Range("L86").Value = 6449
$ws.Range("M86").Value = -3700
$ws.Range("N86").Value = -8695
$ws.Range("H89").Value = 5473.4
$ws.Range("I89").Value = 4823
$ws.Range("J89").Value = 6449
$ws.Range("K89").Value = 24115
$ws.Range("L89").Value = 32245
$ws.Range("M89").Value = -18499
$ws.Range("N89").Value = -43477
$ws.Range("H137").Value = 590570.4399999999
$ws.Range("J137").Value = 2988.125
$ws.Range("L137").Value = 8964.375
$ws.Range("N137").Value = -14064.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3093.3
$ws.Range("I32").Value = 2601.8809
$ws.Range("K32").Value = 2601.8809
$ws.Range("M32").Value = -2314.8809
$ws.Range("H74").Value = 4717.591
$ws.Range("I74").Value = 1386.7273
$ws.Range("J74").Value = 8048.4546
$ws.Range("K74").Value = 1386.7273
$ws.Range("L74").Value = 8048.4546
$ws.Range("M74").Value = -512.7273
$ws.Range("N74").Value = -9796.454600000001
$ws.Range("H77").Value = 4717.591
$ws.Range("I77").Value = 1386.7273
$ws.Range("J77").Value = 8048.4546
$ws.Range("K77").Value = 6933.636500000001
$ws.Range("L77").Value = 40242.273
$ws.Range("M77").Value = -2565.636500000001
$ws.Range("N77").Value = -48978.273
$ws.Range("H97").Value = 2298.6155
$ws.Range("J97").Value = 4032.6
$ws.Range("L97").Value = 4032.6
$ws.Range("N97").Value = -5024.6
$ws.Range("H102").Value = 1839.9546
$ws.Range("I102").Value = 1529.5
$ws.Range("K102").Value = 1529.5
$ws.Range("M102").Value = 92.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1227.7646
$ws.Range("I20").Value = 1015.9
$ws.Range("J20").Value = 1530.4286
$ws.Range("K20").Value = 1015.9
$ws.Range("L20").Value = 1530.4286
$ws.Range("M20").Value = -768.9
$ws.Range("N20").Value = -2024.4286
$ws.Range("H43").Value = 280000
$ws.Range("J43").Value = 280000
$ws.Range("L43").Value = 280000
$ws.Range("N43").Value = -280362
$ws.Range("H86").Value = 9651.929
$ws.Range("I86").Value = 349
$ws.Range("K86").Value = 349
$ws.Range("M86").Value = 774
$ws.Range("H89").Value = 9651.929
$ws.Range("I89").Value = 349
$ws.Range("K89").Value = 1745
$ws.Range("M89").Value = 3871
$ws.Range("H94").Value = 9940.823
$ws.Range("I94").Value = 3598.9092
$ws.Range("K94").Value = 3598.9092
$ws.Range("M94").Value = -3147.9092
$ws.Range("H105").Value = 3636.4285
$ws.Range("I105").Value = 2220.182
$ws.Range("K105").Value = 2220.182
$ws.Range("M105").Value = -473.1819999999998
$ws.Range("H107").Value = 2028.9231
$ws.Range("I107").Value = 1835.9524
$ws.Range("K107").Value = 1835.9524
$ws.Range("M107").Value = 84.0476000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3809.2856
$ws.Range("I2").Value = 6100
$ws.Range("K2").Value = 6100
$ws.Range("M2").Value = -5987
$ws.Range("H11").Value = 2987
$ws.Range("J11").Value = 2987
$ws.Range("L11").Value = 2987
$ws.Range("N11").Value = -3267
$ws.Range("H13").Value = 3999
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("H16").Value = 3227.7058
$ws.Range("I16").Value = 2360.9092
$ws.Range("J16").Value = 4816.8335
$ws.Range("K16").Value = 2360.9092
$ws.Range("L16").Value = 4816.8335
$ws.Range("M16").Value = -2073.9092
$ws.Range("N16").Value = -5390.8335
$ws.Range("H99").Value = 2811.4707
$ws.Range("I99").Value = 2581.7273
$ws.Range("J99").Value = 3232.6667
$ws.Range("K99").Value = 2581.7273
$ws.Range("L99").Value = 3232.6667
$ws.Range("M99").Value = -1083.7273
$ws.Range("N99").Value = -6228.6667
$ws.Range("H107").Value = 941.36365
$ws.Range("I107").Value = 795.1111
$ws.Range("K107").Value = 795.1111
$ws.Range("M107").Value = 1124.8889
$ws.Range("H113").Value = 3227.7058
$ws.Range("I113").Value = 2360.9092
$ws.Range("J113").Value = 4816.8335
$ws.Range("K113").Value = 2360.9092
$ws.Range("L113").Value = 4816.8335
$ws.Range("M113").Value = -190.9092000000001
$ws.Range("N113").Value = -9156.833500000001
$ws.Range("H126").Value = 2811.4707
$ws.Range("I126").Value = 2581.7273
$ws.Range("J126").Value = 3232.6667
$ws.Range("K126").Value = 7745.1819
$ws.Range("L126").Value = 9698.000100000001
$ws.Range("M126").Value = -5275.1819
$ws.Range("N126").Value = -14638.0001
$ws.Range("M13").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 372.27274
$ws.Range("J18").Value = 641.6
$ws.Range("L18").Value = 1924.8
$ws.Range("N18").Value = -2262.8
$ws.Range("H58").Value = 7143.5713
$ws.Range("I58").Value = 4666.6665
$ws.Range("J58").Value = 9001.25
$ws.Range("K58").Value = 13999.9995
$ws.Range("L58").Value = 27003.75
$ws.Range("M58").Value = -13871.9995
$ws.Range("N58").Value = -27259.75
$ws.Range("H98").Value = 26319968
$ws.Range("J98").Value = 31253858
$ws.Range("L98").Value = 93761574
$ws.Range("N98").Value = -93764570
$ws.Range("H109").Value = 7167.8335
$ws.Range("I109").Value = 3501.75
$ws.Range("K109").Value = 10505.25
$ws.Range("M109").Value = -9465.25
$ws.Range("H138").Value = 52653776
$ws.Range("I138").Value = 200000510
$ws.Range("K138").Value = 600001530
$ws.Range("M138").Value = -599996390

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7200.263
$ws.Range("I70").Value = 6886.357
$ws.Range("K70").Value = 6886.357
$ws.Range("M70").Value = -6616.357
$ws.Range("H73").Value = 7200.263
$ws.Range("I73").Value = 6886.357
$ws.Range("K73").Value = 6886.357
$ws.Range("M73").Value = -5950.357
$ws.Range("H80").Value = 1743
$ws.Range("I80").Value = 1250
$ws.Range("K80").Value = 1250
$ws.Range("M80").Value = -252
$ws.Range("H83").Value = 1743
$ws.Range("I83").Value = 1250
$ws.Range("K83").Value = 6250
$ws.Range("M83").Value = -1258
$ws.Range("H97").Value = 1295.4286
$ws.Range("I97").Value = 294
$ws.Range("K97").Value = 294
$ws.Range("M97").Value = 202

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("H93").Value = 2042.3334
$ws.Range("I93").Value = 1955.2727
$ws.Range("K93").Value = 1955.2727
$ws.Range("M93").Value = -707.2727
$ws.Range("H132").Value = 4202.5264
$ws.Range("I132").Value = 3624.5334
$ws.Range("K132").Value = 10873.6002
$ws.Range("M132").Value = -8343.600199999999
$ws.Range("N24").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 36476.8
$ws.Range("J52").Value = 95000
$ws.Range("L52").Value = 95000
$ws.Range("N52").Value = -95452
$ws.Range("H81").Value = 7502
$ws.Range("J81").Value = 4000
$ws.Range("L81").Value = 8000
$ws.Range("N81").Value = -10122
$ws.Range("H84").Value = 7502
$ws.Range("J84").Value = 4000
$ws.Range("L84").Value = 40000
$ws.Range("N84").Value = -50608
$ws.Range("H122").Value = 5887.8887
$ws.Range("I122").Value = 6011.0625
$ws.Range("K122").Value = 18033.1875
$ws.Range("M122").Value = -15583.1875
